$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data rows (2-3) completely first, so stale shared strings
# referenced only by these rows are fully released before the new values are written.
$ws.Range("A2:G3").ClearContents()

# Row 2 - SE151252 / BAP GROUP / Passed
$ws.Range("A2").Value = "SE151252"
$ws.Range("B2").Value = "Kỹ thuật phần mềm"
$ws.Range("C2").Value = "BAP GROUP"
$ws.Range("D2").Value = "Kỹ thuật phần mềm"
$ws.Range("E2").Value = 6.0
$ws.Range("F2").Value = "Nhiệt huyết trong công việc."
$ws.Range("G2").Value = "Passed"

# Row 3 - SE151272 / NASHTECH / Passed
$ws.Range("A3").Value = "SE151272"
$ws.Range("B3").Value = "Kỹ thuật phần mềm"
$ws.Range("C3").Value = "NASHTECH"
$ws.Range("D3").Value = "Kỹ thuật phần mềm"
$ws.Range("E3").Value = 6.0
$ws.Range("F3").Value = "Hoàn thành tốt."
$ws.Range("G3").Value = "Passed"

# Row 4 (new) - SE151262 / MANULIFE / Passed
$ws.Range("A4").Value = "SE151262"
$ws.Range("B4").Value = "Kỹ thuật phần mềm"
$ws.Range("C4").Value = "MANULIFE"
$ws.Range("D4").Value = "Kinh doanh quốc tế"
$ws.Range("E4").Value = 7.0
$ws.Range("F4").Value = "Thực hiện tốt nhiệm vụ giao."
$ws.Range("G4").Value = "Passed"
